$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Cells.Item(2, 4).Value = "23.977.82"
$ws.Cells.Item(2, 5).Value = "  +0.68%  "

$ws.Cells.Item(3, 4).Value = "1.661.10"
$ws.Cells.Item(3, 5).Value = "  +2.86%  "

Set-TextValue $ws.Cells.Item(4, 4) "0.9990"
$ws.Cells.Item(4, 5).Value = "  -0.23%  "

Set-TextValue $ws.Cells.Item(5, 4) "309.47"
$ws.Cells.Item(5, 5).Value = "  +0.94%  "

$ws.Cells.Item(6, 5).Value = "  -0.11%  "

Set-TextValue $ws.Cells.Item(7, 4) "0.3890"
$ws.Cells.Item(7, 5).Value = "  -0.10%  "

Set-TextValue $ws.Cells.Item(8, 4) "0.3839"
$ws.Cells.Item(8, 5).Value = "  +1.12%  "

Set-TextValue $ws.Cells.Item(9, 4) "51.14"
$ws.Cells.Item(9, 5).Value = "  +5.38%  "

Set-TextValue $ws.Cells.Item(10, 4) "1.358"
$ws.Cells.Item(10, 5).Value = "  +0.50%  "

Set-TextValue $ws.Cells.Item(11, 4) "1.002"
$ws.Cells.Item(11, 5).Value = "  +0.05%  "

Set-TextValue $ws.Cells.Item(12, 4) "0.08471"
$ws.Cells.Item(12, 5).Value = "  +0.53%  "

Set-TextValue $ws.Cells.Item(13, 4) "23.98"
$ws.Cells.Item(13, 5).Value = "  +0.97%  "

Set-TextValue $ws.Cells.Item(14, 4) "7.158"
$ws.Cells.Item(14, 5).Value = "  +2.41%  "

Set-TextValue $ws.Cells.Item(15, 4) "7.873"
$ws.Cells.Item(15, 5).Value = "  +6.10%  "

Set-TextValue $ws.Cells.Item(16, 4) "0.00001310"
$ws.Cells.Item(16, 5).Value = "  +3.19%  "

$ws.Cells.Item(17, 4).Value = "1.653.88"
$ws.Cells.Item(17, 5).Value = "  +2.52%  "

Set-TextValue $ws.Cells.Item(18, 4) "94.87"
$ws.Cells.Item(18, 5).Value = "  +1.93%  "

Set-TextValue $ws.Cells.Item(19, 4) "0.07007"
$ws.Cells.Item(19, 5).Value = "  +1.44%  "

Set-TextValue $ws.Cells.Item(20, 4) "19.83"
$ws.Cells.Item(20, 5).Value = "  -0.45%  "

Set-TextValue $ws.Cells.Item(21, 4) "6.911"
$ws.Cells.Item(21, 5).Value = "  +1.89%  "

$ws.Cells.Item(22, 5).Value = "  -0.11%  "

Set-TextValue $ws.Cells.Item(23, 4) "13.62"
$ws.Cells.Item(23, 5).Value = "  +1.85%  "

$ws.Cells.Item(24, 4).Value = "23.956.00"
$ws.Cells.Item(24, 5).Value = "  +0.56%  "

Set-TextValue $ws.Cells.Item(25, 4) "2.492"
$ws.Cells.Item(25, 5).Value = "  +2.68%  "

Set-TextValue $ws.Cells.Item(26, 4) "3.046"
$ws.Cells.Item(26, 5).Value = "  +8.82%  "

Set-TextValue $ws.Cells.Item(27, 4) "22.12"
$ws.Cells.Item(27, 5).Value = "  +0.26%  "

Set-TextValue $ws.Cells.Item(28, 4) "152.67"
$ws.Cells.Item(28, 5).Value = "  -2.79%  "

Set-TextValue $ws.Cells.Item(29, 4) "5.472"
$ws.Cells.Item(29, 5).Value = "  +4.31%  "

Set-TextValue $ws.Cells.Item(30, 4) "139.48"
$ws.Cells.Item(30, 5).Value = "  +0.31%  "

Set-TextValue $ws.Cells.Item(31, 4) "7.804"
$ws.Cells.Item(31, 5).Value = "  +0.99%  "

Set-TextValue $ws.Cells.Item(32, 4) "2.502"
$ws.Cells.Item(32, 5).Value = "  +0.81%  "

$ws.Cells.Item(33, 4).Value = "1.833.10"
$ws.Cells.Item(33, 5).Value = "  +2.36%  "

$ws.Cells.Item(34, 5).Value = "  +7.52%  "

Set-TextValue $ws.Cells.Item(35, 4) "0.08045"
$ws.Cells.Item(35, 5).Value = "  -0.26%  "

Set-TextValue $ws.Cells.Item(36, 4) "0.02966"
$ws.Cells.Item(36, 5).Value = "  +3.64%  "

Set-TextValue $ws.Cells.Item(37, 4) "11.01"
$ws.Cells.Item(37, 5).Value = "  +6.16%  "

Set-TextValue $ws.Cells.Item(38, 4) "6.699"
$ws.Cells.Item(38, 5).Value = "  +2.58%  "

Set-TextValue $ws.Cells.Item(39, 4) "0.2689"
$ws.Cells.Item(39, 5).Value = "  +1.91%  "

Set-TextValue $ws.Cells.Item(40, 4) "0.09124"
$ws.Cells.Item(40, 5).Value = "  -0.02%  "

Set-TextValue $ws.Cells.Item(41, 4) "0.7557"
$ws.Cells.Item(41, 5).Value = "  +1.87%  "

Set-TextValue $ws.Cells.Item(42, 4) "13.53"
$ws.Cells.Item(42, 5).Value = "  +1.40%  "

Set-TextValue $ws.Cells.Item(43, 4) "1.421"

Set-TextValue $ws.Cells.Item(44, 4) "16.26"
$ws.Cells.Item(44, 5).Value = "  +3.36%  "

Set-TextValue $ws.Cells.Item(45, 4) "0.6958"
$ws.Cells.Item(45, 5).Value = "  +2.32%  "

Set-TextValue $ws.Cells.Item(46, 4) "2.467"
$ws.Cells.Item(46, 5).Value = "  +1.54%  "

$ws.Cells.Item(47, 5).Value = "  +0.59%  "

Set-TextValue $ws.Cells.Item(48, 4) "0.9999"
$ws.Cells.Item(48, 5).Value = "  -0.17%  "

Set-TextValue $ws.Cells.Item(49, 4) "0.08293"
$ws.Cells.Item(49, 5).Value = "  +1.01%  "

Set-TextValue $ws.Cells.Item(50, 4) "134.42"
$ws.Cells.Item(50, 5).Value = "  +1.59%  "

Set-TextValue $ws.Cells.Item(51, 4) "1.233"
$ws.Cells.Item(51, 5).Value = "  +4.23%  "
